$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Fill in the "Count" (column E) values for every "On Campus (Residence
#    Halls)" record (rows 29-55) - replicating the Count operation that was
#    already applied to the "On Campus (excluding Residence Halls)" block
#    (rows 2-28).
# ---------------------------------------------------------------------------

$counts = @{
  29 = 197
  30 = 43
  31 = 2
  32 = 20
  33 = 1
  34 = 1
  35 = 0
  36 = $null
  37 = 0
  38 = 167
  39 = 36
  40 = 1
  41 = 23
  42 = 0
  43 = 3
  44 = 1
  45 = $null
  46 = 0
  47 = 212
  48 = 54
  49 = 5
  50 = 35
  51 = 4
  52 = 1
  53 = 0
  54 = $null
  55 = 0
}

# A cell that already carries the same number style used throughout column E
# (rows 2-28) - used as the format donor so the newly written cells pick up
# the identical style instead of the plain "no style" default.
$formatDonor = $ws.Range("E2")

foreach ($row in 29..55) {
    $cell = $ws.Range("E$row")
    $value = $counts[$row]

    if ($null -ne $value) {
        $cell.Value2 = $value
    }

    $formatDonor.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Update the view/selection state: the sheet was scrolled down a bit
#    further and a different cell ended up active/selected.
# ---------------------------------------------------------------------------

$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1
[void]$ws.Range("F48").Select()

$win.Left = 9340
$win.Top = 1380
